$wb = $excel.ActiveWorkbook

function Set-PctText {
    param($ws, $cellRef, $text)
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
}

# Sheet: Summary (all 15 data rows)
$ws = $wb.Worksheets.Item("Summary")
Set-PctText $ws "G2" "+236.65%"
$ws.Range("M2").Value = 4
Set-PctText $ws "G3" "+218.11%"
$ws.Range("M3").Value = 4
Set-PctText $ws "G4" "+99.85%"
$ws.Range("M4").Value = 4
Set-PctText $ws "G5" "+29.33%"
$ws.Range("M5").Value = 4
Set-PctText $ws "G6" "+3.12%"
$ws.Range("M6").Value = 4
Set-PctText $ws "G7" "+17.67%"
$ws.Range("M7").Value = 4
Set-PctText $ws "G8" "+29.45%"
$ws.Range("M8").Value = 4
Set-PctText $ws "G9" "+14.60%"
$ws.Range("M9").Value = 4
Set-PctText $ws "G10" "+5.46%"
$ws.Range("M10").Value = 4
Set-PctText $ws "G11" "+2.03%"
$ws.Range("M11").Value = 4
Set-PctText $ws "G12" "+11.50%"
$ws.Range("M12").Value = 4
Set-PctText $ws "G13" "+118.08%"
$ws.Range("M13").Value = 4
Set-PctText $ws "G14" "+165.01%"
$ws.Range("M14").Value = 4
Set-PctText $ws "G15" "+4.46%"
$ws.Range("M15").Value = 4
Set-PctText $ws "G16" "+36.30%"
$ws.Range("M16").Value = 4

# Sheet: Pattern1-Pure Data (rows 2-6 mirror Summary rows 2-6)
$ws = $wb.Worksheets.Item("Pattern1-Pure Data")
Set-PctText $ws "G2" "+236.65%"
$ws.Range("M2").Value = 4
Set-PctText $ws "G3" "+218.11%"
$ws.Range("M3").Value = 4
Set-PctText $ws "G4" "+99.85%"
$ws.Range("M4").Value = 4
Set-PctText $ws "G5" "+29.33%"
$ws.Range("M5").Value = 4
Set-PctText $ws "G6" "+3.12%"
$ws.Range("M6").Value = 4

# Sheet: Pattern2-Data+Technical (rows 2-6 mirror Summary rows 7-11)
$ws = $wb.Worksheets.Item("Pattern2-Data+Technical")
Set-PctText $ws "G2" "+17.67%"
$ws.Range("M2").Value = 4
Set-PctText $ws "G3" "+29.45%"
$ws.Range("M3").Value = 4
Set-PctText $ws "G4" "+14.60%"
$ws.Range("M4").Value = 4
Set-PctText $ws "G5" "+5.46%"
$ws.Range("M5").Value = 4
Set-PctText $ws "G6" "+2.03%"
$ws.Range("M6").Value = 4

# Sheet: Pattern3-Data+News (rows 2-6 mirror Summary rows 12-16)
$ws = $wb.Worksheets.Item("Pattern3-Data+News")
Set-PctText $ws "G2" "+11.50%"
$ws.Range("M2").Value = 4
Set-PctText $ws "G3" "+118.08%"
$ws.Range("M3").Value = 4
Set-PctText $ws "G4" "+165.01%"
$ws.Range("M4").Value = 4
Set-PctText $ws "G5" "+4.46%"
$ws.Range("M5").Value = 4
Set-PctText $ws "G6" "+36.30%"
$ws.Range("M6").Value = 4
